$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 88988.8
$ws.Range("I6").Value = 500.33334
$ws.Range("J6").Value = 221721.5
$ws.Range("K6").Value = 1501.00002
$ws.Range("L6").Value = 665164.5
$ws.Range("M6").Value = -1389.00002
$ws.Range("N6").Value = -665388.5
$ws.Range("H12").Value = 450.82144
$ws.Range("I12").Value = 420.14285
$ws.Range("J12").Value = 542.8570999999999
$ws.Range("K12").Value = 420.14285
$ws.Range("L12").Value = 542.8570999999999
$ws.Range("M12").Value = -250.14285
$ws.Range("N12").Value = -882.8570999999999
$ws.Range("H21").Value = 42500
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 80000
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 80000
$ws.Range("M21").Value = -4532
$ws.Range("N21").Value = -80936
$ws.Range("H23").Value = 42500
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 80000
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 80000
$ws.Range("M23").Value = -4766
$ws.Range("N23").Value = -80468
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 100
$ws.Range("K29").Value = 300
$ws.Range("M29").Value = -19
$ws.Range("H45").Value = 2933.3333
$ws.Range("J45").Value = 2933.3333
$ws.Range("L45").Value = 8799.999899999999
$ws.Range("N45").Value = -9183.999899999999
$ws.Range("H92").Value = 1075
$ws.Range("I92").Value = 599.1667
$ws.Range("J92").Value = 2502.5
$ws.Range("K92").Value = 599.1667
$ws.Range("L92").Value = 2502.5
$ws.Range("M92").Value = 648.8333
$ws.Range("N92").Value = -4998.5
$ws.Range("H129").Value = 916.24243
$ws.Range("J129").Value = 1050.7778
$ws.Range("L129").Value = 3152.3334
$ws.Range("N129").Value = -13152.3334
$ws.Range("H132").Value = 944787.7
$ws.Range("I132").Value = 2718.3713
$ws.Range("J132").Value = 2884342
$ws.Range("K132").Value = 8155.113899999999
$ws.Range("L132").Value = 8653026
$ws.Range("M132").Value = -5625.113899999999
$ws.Range("N132").Value = -8658086
$ws.Range("H136").Value = 50780
$ws.Range("J136").Value = 50780
$ws.Range("L136").Value = 50780
$ws.Range("N136").Value = -60980
$ws.Range("H137").Value = 2175478.5
$ws.Range("I137").Value = 2942262
$ws.Range("J137").Value = 2925.4167
$ws.Range("K137").Value = 8826786
$ws.Range("L137").Value = 8776.250100000001
$ws.Range("M137").Value = -8824236
$ws.Range("N137").Value = -13876.2501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 22267920
$ws.Range("I61").Value = 27055168
$ws.Range("J61").Value = 126901.75
$ws.Range("K61").Value = 27055168
$ws.Range("L61").Value = 126901.75
$ws.Range("M61").Value = -27054956
$ws.Range("N61").Value = -127325.75
$ws.Range("H74").Value = 10480820
$ws.Range("I74").Value = 13374655
$ws.Range("J74").Value = 145694.86
$ws.Range("K74").Value = 13374655
$ws.Range("L74").Value = 145694.86
$ws.Range("M74").Value = -13373781
$ws.Range("N74").Value = -147442.86
$ws.Range("H77").Value = 10480820
$ws.Range("I77").Value = 13374655
$ws.Range("J77").Value = 145694.86
$ws.Range("K77").Value = 66873275
$ws.Range("L77").Value = 728474.2999999999
$ws.Range("M77").Value = -66868907
$ws.Range("N77").Value = -737210.2999999999
$ws.Range("H132").Value = 62697.617
$ws.Range("I132").Value = 50894.25
$ws.Range("J132").Value = 79559.57000000001
$ws.Range("K132").Value = 152682.75
$ws.Range("L132").Value = 238678.71
$ws.Range("M132").Value = -150152.75
$ws.Range("N132").Value = -243738.71
$ws.Range("H134").Value = 49429
$ws.Range("J134").Value = 49429
$ws.Range("L134").Value = 49429
$ws.Range("N134").Value = -59569
$ws.Range("H136").Value = 22267920
$ws.Range("I136").Value = 27055168
$ws.Range("J136").Value = 126901.75
$ws.Range("K136").Value = 81165504
$ws.Range("L136").Value = 380705.25
$ws.Range("M136").Value = -81162954
$ws.Range("N136").Value = -385805.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 915.86664
$ws.Range("I99").Value = 845.7
$ws.Range("J99").Value = 1056.2
$ws.Range("K99").Value = 845.7
$ws.Range("L99").Value = 1056.2
$ws.Range("M99").Value = 652.3
$ws.Range("N99").Value = -4052.2
$ws.Range("H132").Value = 37000
$ws.Range("I132").Value = 37000
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 37000
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -31940
$ws.Range("N132").Value = $null
$ws.Range("H138").Value = 31631.666
$ws.Range("J138").Value = 31631.666
$ws.Range("L138").Value = 31631.666
$ws.Range("N138").Value = -41911.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3235.4443
$ws.Range("I31").Value = 1577.375
$ws.Range("J31").Value = 16500
$ws.Range("K31").Value = 1577.375
$ws.Range("L31").Value = 16500
$ws.Range("M31").Value = -1282.375
$ws.Range("N31").Value = -17090
$ws.Range("H34").Value = 3235.4443
$ws.Range("I34").Value = 1577.375
$ws.Range("J34").Value = 16500
$ws.Range("K34").Value = 1577.375
$ws.Range("L34").Value = 16500
$ws.Range("M34").Value = -1375.375
$ws.Range("N34").Value = -16904
$ws.Range("H58").Value = 34485364
$ws.Range("I58").Value = 66669908
$ws.Range("J58").Value = 1923.1428
$ws.Range("K58").Value = 66669908
$ws.Range("L58").Value = 1923.1428
$ws.Range("M58").Value = -66669705
$ws.Range("N58").Value = -2329.1428
$ws.Range("H132").Value = 31731.03
$ws.Range("I132").Value = 2333.842
$ws.Range("J132").Value = 68967.47
$ws.Range("K132").Value = 7001.526
$ws.Range("L132").Value = 206902.41
$ws.Range("M132").Value = -4471.526
$ws.Range("N132").Value = -211962.41
$ws.Range("H136").Value = 34485364
$ws.Range("I136").Value = 66669908
$ws.Range("J136").Value = 1923.1428
$ws.Range("K136").Value = 200009724
$ws.Range("L136").Value = 5769.428400000001
$ws.Range("M136").Value = -200007174
$ws.Range("N136").Value = -10869.4284

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5415.5625
$ws.Range("I3").Value = 5554.0835
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 16662.2505
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = -16550.2505
$ws.Range("N3").Value = -15224
$ws.Range("H113").Value = 445.9
$ws.Range("I113").Value = 337.8
$ws.Range("K113").Value = 1013.4
$ws.Range("M113").Value = 1156.6
$ws.Range("H131").Value = 1020.25714
$ws.Range("I131").Value = 403.33334
$ws.Range("J131").Value = 1111.2787
$ws.Range("K131").Value = 1210.00002
$ws.Range("L131").Value = 3333.8361
$ws.Range("M131").Value = 3829.99998
$ws.Range("N131").Value = -13413.8361

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 37136.332
$ws.Range("J141").Value = 37136.332
$ws.Range("L141").Value = 37136.332
$ws.Range("N141").Value = -47496.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3063.7856
$ws.Range("I122").Value = 2786.6316
$ws.Range("J122").Value = 3648.889
$ws.Range("K122").Value = 8359.8948
$ws.Range("L122").Value = 10946.667
$ws.Range("M122").Value = -5909.8948
$ws.Range("N122").Value = -15846.667
$ws.Range("H132").Value = 23522.412
$ws.Range("I132").Value = 1491.7297
$ws.Range("J132").Value = 114093
$ws.Range("K132").Value = 4475.189100000001
$ws.Range("L132").Value = 342279
$ws.Range("M132").Value = -1945.189100000001
$ws.Range("N132").Value = -347339
$ws.Range("H136").Value = 68291.92999999999
$ws.Range("I136").Value = 46671.047
$ws.Range("J136").Value = 127749.375
$ws.Range("K136").Value = 140013.141
$ws.Range("L136").Value = 383248.125
$ws.Range("M136").Value = -137463.141
$ws.Range("N136").Value = -388348.125
$ws.Range("H138").Value = 51519.332
$ws.Range("J138").Value = 51519.332
$ws.Range("L138").Value = 51519.332
$ws.Range("N138").Value = -61799.332
$ws.Range("H140").Value = 43232.25
$ws.Range("J140").Value = 43232.25
$ws.Range("L140").Value = 43232.25
$ws.Range("N140").Value = -53592.25
$ws.Range("H141").Value = 59540.91
$ws.Range("J141").Value = 59540.91
$ws.Range("L141").Value = 59540.91
$ws.Range("N141").Value = -69900.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 43657.25
$ws.Range("I46").Value = 18000
$ws.Range("J46").Value = 47322.57
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 47322.57
$ws.Range("M46").Value = -17769
$ws.Range("N46").Value = -47784.57
$ws.Range("H100").Value = 63546.5
$ws.Range("I100").Value = 62917.625
$ws.Range("J100").Value = 64175.375
$ws.Range("K100").Value = 125835.25
$ws.Range("L100").Value = 128350.75
$ws.Range("M100").Value = -125294.25
$ws.Range("N100").Value = -129432.75
$ws.Range("H132").Value = 47005.043
$ws.Range("I132").Value = 38999.207
$ws.Range("J132").Value = 60662.06
$ws.Range("K132").Value = 116997.621
$ws.Range("L132").Value = 181986.18
$ws.Range("M132").Value = -114467.621
$ws.Range("N132").Value = -187046.18
$ws.Range("H134").Value = 43657.25
$ws.Range("I134").Value = 18000
$ws.Range("J134").Value = 47322.57
$ws.Range("K134").Value = 54000
$ws.Range("L134").Value = 141967.71
$ws.Range("M134").Value = -51465
$ws.Range("N134").Value = -147037.71
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws.Range("H136").Value = 37772.445
$ws.Range("I136").Value = 31398.697
$ws.Range("J136").Value = 46917.39
$ws.Range("K136").Value = 94196.091
$ws.Range("L136").Value = 140752.17
$ws.Range("M136").Value = -91646.091
$ws.Range("N136").Value = -145852.17
$ws.Range("H137").Value = 57611.316
$ws.Range("J137").Value = 57611.316
$ws.Range("L137").Value = 57611.316
$ws.Range("N137").Value = -67811.31599999999
$ws.Range("H140").Value = 55821.145
$ws.Range("J140").Value = 55821.145
$ws.Range("L140").Value = 55821.145
$ws.Range("N140").Value = -66181.14499999999
$ws.Range("H141").Value = 65367.5
$ws.Range("J141").Value = 65367.5
$ws.Range("L141").Value = 65367.5
$ws.Range("N141").Value = -75727.5
